$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 377.14285
$ws.Range("J19").Value = 255.4
$ws.Range("L19").Value = 255.4
$ws.Range("N19").Value = -605.4
$ws.Range("H33").Value = 1009.4545
$ws.Range("I33").Value = 973.7778
$ws.Range("K33").Value = 973.7778
$ws.Range("M33").Value = -744.7778
$ws.Range("H40").Value = 4953.375
$ws.Range("I40").Value = 3083.3333
$ws.Range("J40").Value = 6075.4
$ws.Range("K40").Value = 3083.3333
$ws.Range("L40").Value = 6075.4
$ws.Range("M40").Value = -2908.3333
$ws.Range("N40").Value = -6425.4
$ws.Range("H86").Value = 3763842.2
$ws.Range("I86").Value = 4298.5557
$ws.Range("J86").Value = 10531021
$ws.Range("K86").Value = 4298.5557
$ws.Range("L86").Value = 10531021
$ws.Range("M86").Value = -3175.5557
$ws.Range("N86").Value = -10533267
$ws.Range("H89").Value = 3763842.2
$ws.Range("I89").Value = 4298.5557
$ws.Range("J89").Value = 10531021
$ws.Range("K89").Value = 21492.7785
$ws.Range("L89").Value = 52655105
$ws.Range("M89").Value = -15876.7785
$ws.Range("N89").Value = -52666337
$ws.Range("H127").Value = 11968.808
$ws.Range("I127").Value = 1775
$ws.Range("K127").Value = 5325
$ws.Range("M127").Value = -365

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3949.3333
$ws.Range("I61").Value = 2621.0625
$ws.Range("K61").Value = 2621.0625
$ws.Range("M61").Value = -2409.0625
$ws.Range("H74").Value = 2773.0605
$ws.Range("I74").Value = 2589.8215
$ws.Range("K74").Value = 2589.8215
$ws.Range("M74").Value = -1715.8215
$ws.Range("H77").Value = 2773.0605
$ws.Range("I77").Value = 2589.8215
$ws.Range("K77").Value = 12949.1075
$ws.Range("M77").Value = -8581.1075
$ws.Range("H122").Value = 4808.6177
$ws.Range("I122").Value = 3965.9285
$ws.Range("K122").Value = 11897.7855
$ws.Range("M122").Value = -9447.7855
$ws.Range("H136").Value = 3949.3333
$ws.Range("I136").Value = 2621.0625
$ws.Range("K136").Value = 7863.1875
$ws.Range("M136").Value = -5313.1875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 261.3158
$ws.Range("I7").Value = 54.8
$ws.Range("J7").Value = 335.07144
$ws.Range("K7").Value = 54.8
$ws.Range("L7").Value = 335.07144
$ws.Range("M7").Value = 58.2
$ws.Range("N7").Value = -561.0714399999999
$ws.Range("H58").Value = 2691
$ws.Range("I58").Value = 2189.6
$ws.Range("J58").Value = 3944.5
$ws.Range("K58").Value = 2189.6
$ws.Range("L58").Value = 3944.5
$ws.Range("M58").Value = -1986.6
$ws.Range("N58").Value = -4350.5
$ws.Range("H132").Value = 2790.7778
$ws.Range("I132").Value = 2446.0715
$ws.Range("K132").Value = 7338.2145
$ws.Range("M132").Value = -4808.2145
$ws.Range("H136").Value = 2691
$ws.Range("I136").Value = 2189.6
$ws.Range("J136").Value = 3944.5
$ws.Range("K136").Value = 6568.799999999999
$ws.Range("L136").Value = 11833.5
$ws.Range("M136").Value = -4018.799999999999
$ws.Range("N136").Value = -16933.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 45718650
$ws.Range("I4").Value = 52937228
$ws.Range("K4").Value = 158811684
$ws.Range("M4").Value = -158811572
$ws.Range("H57").Value = 61299.715
$ws.Range("J57").Value = 61299.715
$ws.Range("L57").Value = 183899.145
$ws.Range("N57").Value = -185017.145
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").Value = $null
$ws.Range("H117").Value = 329.8125
$ws.Range("I117").Value = 141.16667
$ws.Range("J117").Value = 443
$ws.Range("K117").Value = 423.50001
$ws.Range("L117").Value = 1329
$ws.Range("M117").Value = 3018.49999
$ws.Range("N117").Value = -8213
$ws.Range("H118").Value = 3278.111
$ws.Range("I118").Value = 300.6
$ws.Range("K118").Value = 901.8000000000001
$ws.Range("M118").Value = 341.1999999999999
$ws.Range("H121").Value = 1908.3529
$ws.Range("I121").Value = 252.33333
$ws.Range("J121").Value = 2263.2144
$ws.Range("K121").Value = 756.99999
$ws.Range("L121").Value = 6789.6432
$ws.Range("M121").Value = 553.00001
$ws.Range("N121").Value = -9409.643199999999
$ws.Range("H129").Value = 2212.8462
$ws.Range("I129").Value = 1259.6666
$ws.Range("J129").Value = 2498.8
$ws.Range("K129").Value = 3778.9998
$ws.Range("L129").Value = 7496.400000000001
$ws.Range("M129").Value = 1221.0002
$ws.Range("N129").Value = -17496.4
$ws.Range("H131").Value = 6131.3335
$ws.Range("J131").Value = 9766.666999999999
$ws.Range("L131").Value = 29300.001
$ws.Range("N131").Value = -39380.001
$ws.Range("H132").Value = 1150
$ws.Range("J132").Value = 1200
$ws.Range("L132").Value = 10800
$ws.Range("N132").Value = -15860

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20905278
$ws.Range("I80").Value = 224617
$ws.Range("K80").Value = 224617
$ws.Range("M80").Value = -223619
$ws.Range("H83").Value = 20905278
$ws.Range("I83").Value = 224617
$ws.Range("K83").Value = 1123085
$ws.Range("M83").Value = -1118093
$ws.Range("H104").Value = 60000
$ws.Range("J104").Value = 60000
$ws.Range("L104").Value = 60000
$ws.Range("N104").Value = -66988
$ws.Range("H108").Value = 50684
$ws.Range("J108").Value = 50684
$ws.Range("L108").Value = 50684
$ws.Range("N108").Value = -58364
$ws.Range("H122").Value = 4568.9165
$ws.Range("I122").Value = 3536.182
$ws.Range("J122").Value = 5442.769
$ws.Range("K122").Value = 10608.546
$ws.Range("L122").Value = 16328.307
$ws.Range("M122").Value = -8158.545999999998
$ws.Range("N122").Value = -21228.307
$ws.Range("H132").Value = 3482.3428
$ws.Range("I132").Value = 2627.261
$ws.Range("J132").Value = 5121.25
$ws.Range("K132").Value = 7881.782999999999
$ws.Range("L132").Value = 15363.75
$ws.Range("M132").Value = -5351.782999999999
$ws.Range("N132").Value = -20423.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6030.6
$ws.Range("I7").Value = 7710.4
$ws.Range("J7").Value = 2671
$ws.Range("K7").Value = 7710.4
$ws.Range("L7").Value = 2671
$ws.Range("M7").Value = -7598.4
$ws.Range("N7").Value = -2895
$ws.Range("H55").Value = 1063.1724
$ws.Range("J55").Value = 558.06665
$ws.Range("L55").Value = 558.06665
$ws.Range("N55").Value = -904.06665
$ws.Range("H93").Value = 314899.47
$ws.Range("I93").Value = 2513.7222
$ws.Range("K93").Value = 2513.7222
$ws.Range("M93").Value = -1265.7222
$ws.Range("H124").Value = 500000
$ws.Range("J124").Value = 500000
$ws.Range("L124").Value = 500000
$ws.Range("N124").Value = -509820
$ws.Range("H126").Value = 6030.6
$ws.Range("I126").Value = 7710.4
$ws.Range("J126").Value = 2671
$ws.Range("K126").Value = 23131.2
$ws.Range("L126").Value = 8013
$ws.Range("M126").Value = -20661.2
$ws.Range("N126").Value = -12953
$ws.Range("H132").Value = 4271.148
$ws.Range("I132").Value = 3438.4707
$ws.Range("K132").Value = 10315.4121
$ws.Range("M132").Value = -7785.4121
$ws.Range("H136").Value = 5893.4136
$ws.Range("I136").Value = 4508.3687
$ws.Range("J136").Value = 8525
$ws.Range("K136").Value = 13525.1061
$ws.Range("L136").Value = 25575
$ws.Range("M136").Value = -10975.1061
$ws.Range("N136").Value = -30675

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H60").Value = 19999
$ws.Range("J60").Value = 19999
$ws.Range("L60").Value = 19999
$ws.Range("N60").Value = -21643
$ws.Range("H81").Value = 23817802
$ws.Range("I81").Value = 11153
$ws.Range("K81").Value = 22306
$ws.Range("M81").Value = -21245
$ws.Range("H84").Value = 23817802
$ws.Range("I84").Value = 11153
$ws.Range("K84").Value = 111530
$ws.Range("M84").Value = -106226
$ws.Range("H86").Value = 83871.5
$ws.Range("J86").Value = 83871.5
$ws.Range("L86").Value = 83871.5
$ws.Range("N86").Value = -86117.5
$ws.Range("H89").Value = 83871.5
$ws.Range("J89").Value = 83871.5
$ws.Range("L89").Value = 419357.5
$ws.Range("N89").Value = -430589.5
$ws.Range("H126").Value = 2683.8333
$ws.Range("I126").Value = 2626.7273
$ws.Range("K126").Value = 7880.1819
$ws.Range("M126").Value = -5410.1819
$ws.Range("H132").Value = 2647.1052
$ws.Range("I132").Value = 2114.9
$ws.Range("J132").Value = 3238.4443
$ws.Range("K132").Value = 6344.700000000001
$ws.Range("L132").Value = 9715.332900000001
$ws.Range("M132").Value = -3814.700000000001
$ws.Range("N132").Value = -14775.3329
